$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Exiobase")
$ws2 = $wb.Worksheets.Item("Deutsch")
$ws3 = $wb.Worksheets.Item("English")

# --- Section: Selection / Auswahl ---
$ws1.Range("A13").Value = "Selection"
$ws1.Range("B13").Value = "Selection"
$ws2.Range("A13").Value = "Selection"
$ws2.Range("B13").Value = "Auswahl"

$ws1.Range("A14").Value = "Visualisation"
$ws1.Range("A15").Value = "Settings"
$ws1.Range("B14").Value = "Visualisation"
$ws1.Range("B15").Value = "Settings"
$ws2.Range("A14").Value = "Visualisation"
$ws2.Range("A15").Value = "Settings"
$ws2.Range("B14").Value = "Visualisierung"
$ws2.Range("B15").Value = "Einstellungen"

# --- Section: General Settings / Grundeinstellungen ---
$ws1.Range("A16").Value = "General Settings"
$ws1.Range("B16").Value = "General Settings"
$ws2.Range("A16").Value = "General Settings"
$ws2.Range("B16").Value = "Grundeinstellungen"

$ws1.Range("A17").Value = "Language"
$ws1.Range("A18").Value = "Year"
$ws1.Range("B17").Value = "Language"
$ws1.Range("B18").Value = "Year"
$ws2.Range("A17").Value = "Language"
$ws2.Range("A18").Value = "Year"
$ws2.Range("B17").Value = "Sprache"
$ws2.Range("B18").Value = "Jahr"

$ws1.Range("A19").Value = "Show Indices"
$ws1.Range("B19").Value = "Show Indices"
$ws2.Range("A19").Value = "Show Indices"
$ws2.Range("B19").Value = "Indices anzeigen"

# --- Sheet3 (English) mirrors Sheet1 content for rows 13-19 ---
$ws3.Range("A13").Value = "Selection"
$ws3.Range("B13").Value = "Selection"
$ws3.Range("A14").Value = "Visualisation"
$ws3.Range("A15").Value = "Settings"
$ws3.Range("B14").Value = "Visualisation"
$ws3.Range("B15").Value = "Settings"
$ws3.Range("A16").Value = "General Settings"
$ws3.Range("B16").Value = "General Settings"
$ws3.Range("A17").Value = "Language"
$ws3.Range("A18").Value = "Year"
$ws3.Range("B17").Value = "Language"
$ws3.Range("B18").Value = "Year"
$ws3.Range("A19").Value = "Show Indices"
$ws3.Range("B19").Value = "Show Indices"

# --- Sheet3 extra duplicated block rows 20-23 (mirrors Sheet2's B column for rows 16-19) ---
$ws3.Range("A20").Value = "General Settings"
$ws3.Range("B20").Value = "Grundeinstellungen"
$ws3.Range("A21").Value = "Language"
$ws3.Range("B21").Value = "Sprache"
$ws3.Range("A22").Value = "Year"
$ws3.Range("B22").Value = "Jahr"
$ws3.Range("A23").Value = "Show Indices"
$ws3.Range("B23").Value = "Indices anzeigen"

Write-Host "done values"
